$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '29.409.31'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.847.56'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '240.70'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '0.6288'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.07692'
$ws.Range("E8").Value = '  +2.33%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.2921'
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '24.79'
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07739'
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.847.68'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '5.028'
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '0.00001076'
$ws.Range("E14").Value = '  +3.09%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.6799'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '83.59'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '6.174'
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '29.448.30'
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '228.21'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '7.408'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '157.24'
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '17.70'
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '1.356'
$ws.Range("E28").Value = '  +6.32%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '1.464'
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '0.05688'
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '4.116'
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '4.028'
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.7072'
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '1.225.17'
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '6.536'
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.9068'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '101.72'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '66.14'
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '7.165'
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.4016'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '9.034'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.1146'
$ws.Range("E49").Value = '  +2.48%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '1.673'
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.05713'
$ws.Range("E51").Value = '  +0.12%  '
